$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B..E values in row 1 (subject/condition codes)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) updated values
$ws.Range("B2").Value = 47.597757901346469
$ws.Range("C2").Value = 55.718370135679095
$ws.Range("D2").Value = 50.557737295124724
$ws.Range("E2").Value = 56.125127876082288

# Row 3 (STR) updated values
$ws.Range("B3").Value = 44.981597518379402
$ws.Range("C3").Value = 45.876060714702433
$ws.Range("D3").Value = 45.291664428454979
$ws.Range("E3").Value = 55.362706001604877

# Selection changed to B1:E3
[void]$ws.Range("B1:E3").Select()
